# Change 1: remove the bold "Learning Targets:" label and the following space,
# leaving "(All evidence ...)" as the start of the paragraph (unformatted).
$d = $word.ActiveDocument
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Learning Targets: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Output "Change1 found: $found1"

# Change 2: rewrite the tail of the long reflection paragraph.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("I have demonstrated proficiency for the majority of the learning targets that we have learned so far. I have revised code from my labs and challenges when the opportunity was given. Additionally when I did not have to submit a revision, I still improved my code from my peer reviews and comments from Dr. Theobold. My revisions have helped me not make the same mistake in the next labs and challenges. I have extended my thinking to all the Challenge assignments so far, and have taken risks by completing tasks that I did not know would be correct. I communicate with my team and Dr. Theobold, asking questions when there are confusions, and asking for additional help in office hours and in the Discord. I have been present, respectful, prepared for class and completed most of the assignments. Lastly, I have given respectful peer code reviews and put my best effort to provide praise and suggestions for improvements. ", $true, $false, $false, $false, $false, $true, 1, $false, "I have demonstrated proficiency for the majority of the learning targets that we have learned so far, revised code from my labs and challenges when the opportunity was given, have extended my thinking to all the Challenge assignments, and have taken risks by completing tasks that I did not know would be correct. I communicate with my team and Dr. Theobold, asking questions when there are confusions, and asking for additional help in office hours and in the Discord. I have been present, respectful, prepared for class and completed all of the assignments. ", 2)
Write-Output "Change2 found: $found2"
